{"js": "// Replace the date line and the 25 \"a\u00f7b=\" division prompts in the table\n// with their updated values, per the commit's diff. Every old value is\n// unique within the document and none of the new values collide with an\n// old value, so a plain ordered list of exact, case-sensitive\n// search-and-replace operations is sufficient and order-independent.\nconst replacements = [\n  [\"2023-09-20 Wednesday\", \"2023-09-21 Thursday\"],\n  [\"60\u00f79=\", \"41\u00f79=\"],\n  [\"58\u00f73=\", \"16\u00f73=\"],\n  [\"28\u00f74=\", \"45\u00f74=\"],\n  [\"54\u00f72=\", \"69\u00f78=\"],\n  [\"46\u00f72=\", \"69\u00f74=\"],\n  [\"86\u00f75=\", \"17\u00f79=\"],\n  [\"59\u00f79=\", \"86\u00f78=\"],\n  [\"94\u00f76=\", \"19\u00f72=\"],\n  [\"12\u00f75=\", \"82\u00f73=\"],\n  [\"41\u00f78=\", \"37\u00f74=\"],\n  [\"33\u00f74=\", \"72\u00f75=\"],\n  [\"56\u00f76=\", \"52\u00f77=\"],\n  [\"15\u00f78=\", \"97\u00f79=\"],\n  [\"52\u00f78=\", \"13\u00f72=\"],\n  [\"23\u00f74=\", \"50\u00f79=\"],\n  [\"84\u00f77=\", \"21\u00f72=\"],\n  [\"18\u00f72=\", \"87\u00f74=\"],\n  [\"62\u00f76=\", \"56\u00f75=\"],\n  [\"32\u00f74=\", \"12\u00f72=\"],\n  [\"86\u00f73=\", \"47\u00f77=\"],\n  [\"44\u00f75=\", \"75\u00f73=\"],\n  [\"44\u00f72=\", \"34\u00f74=\"],\n  [\"66\u00f73=\", \"95\u00f77=\"],\n  [\"30\u00f78=\", \"55\u00f77=\"],\n  [\"42\u00f72=\", \"30\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 \"a\u00f7b=\" division prompts in the table\n# with their updated values, per the commit's diff. Every old value is\n# unique within the document and none of the new values collide with an\n# old value, so a plain ordered list of exact, case-sensitive\n# find/replace-all operations is sufficient and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-20 Wednesday\", \"2023-09-21 Thursday\"),\n    @(\"60\u00f79=\", \"41\u00f79=\"),\n    @(\"58\u00f73=\", \"16\u00f73=\"),\n    @(\"28\u00f74=\", \"45\u00f74=\"),\n    @(\"54\u00f72=\", \"69\u00f78=\"),\n    @(\"46\u00f72=\", \"69\u00f74=\"),\n    @(\"86\u00f75=\", \"17\u00f79=\"),\n    @(\"59\u00f79=\", \"86\u00f78=\"),\n    @(\"94\u00f76=\", \"19\u00f72=\"),\n    @(\"12\u00f75=\", \"82\u00f73=\"),\n    @(\"41\u00f78=\", \"37\u00f74=\"),\n    @(\"33\u00f74=\", \"72\u00f75=\"),\n    @(\"56\u00f76=\", \"52\u00f77=\"),\n    @(\"15\u00f78=\", \"97\u00f79=\"),\n    @(\"52\u00f78=\", \"13\u00f72=\"),\n    @(\"23\u00f74=\", \"50\u00f79=\"),\n    @(\"84\u00f77=\", \"21\u00f72=\"),\n    @(\"18\u00f72=\", \"87\u00f74=\"),\n    @(\"62\u00f76=\", \"56\u00f75=\"),\n    @(\"32\u00f74=\", \"12\u00f72=\"),\n    @(\"86\u00f73=\", \"47\u00f77=\"),\n    @(\"44\u00f75=\", \"75\u00f73=\"),\n    @(\"44\u00f72=\", \"34\u00f74=\"),\n    @(\"66\u00f73=\", \"95\u00f77=\"),\n    @(\"30\u00f78=\", \"55\u00f77=\"),\n    @(\"42\u00f72=\", \"30\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
